$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new shared strings in the same order they were first introduced
# (matches the authoring order reflected in xl/sharedStrings.xml: ExitRoom,
# EvilTwinsRoom, PPRoom, BossRoom appended after the pre-existing
# StartingRoom/EmptyWoodsPath entries).
$ws.Range("D1").Value = "ExitRoom"
$ws.Range("D3").Value = "EvilTwinsRoom"
$ws.Range("C2").Value = "PPRoom"
$ws.Range("D2").Value = "BossRoom"

# Fill in the remaining cells of the redesigned map
$ws.Range("A3").Value = "PPRoom"
$ws.Range("C3").Value = "EmptyWoodsPath"
$ws.Range("E3").Value = "EvilTwinsRoom"

$ws.Range("A4").Value = "EmptyWoodsPath"
$ws.Range("C4").Value = "EvilTwinsRoom"
$ws.Range("E4").Value = "PPRoom"

$ws.Range("A5").Value = "EvilTwinsRoom"
$ws.Range("B5").Value = "EvilTwinsRoom"
$ws.Range("C5").Value = "StartingRoom"
$ws.Range("D5").Value = "EmptyWoodsPath"
$ws.Range("E5").Value = "EmptyWoodsPath"

$ws.Range("A6").Value = "PPRoom"
$ws.Range("C6").Value = "PPRoom"
$ws.Range("E6").Value = "EvilTwinsRoom"

$ws.Range("A7").Value = "EmptyWoodsPath"
$ws.Range("C7").Value = "PPRoom"
$ws.Range("D7").Value = "EmptyWoodsPath"
$ws.Range("E7").Value = "EvilTwinsRoom"

$ws.Range("B8").Value = "PPRoom"
$ws.Range("C8").Value = "EmptyWoodsPath"

# Update selection to D2
$ws.Range("D2").Select()
